# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 97 (Florida King, Tercera, 14-kilo box,
# Provincia de Limarí), shifting the existing rows 97-169 down to 98-170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(97).Insert()

$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 44512
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100103
$ws.Range("H97").Value = "Frutos de hueso (carozo)"
$ws.Range("I97").Value = 100103004
$ws.Range("J97").Value = "Durazno"
$ws.Range("K97").Value = "Florida King"
$ws.Range("L97").Value = "Tercera"
$ws.Range("M97").Value = 300
$ws.Range("N97").Value = 18000
$ws.Range("O97").Value = 18000
$ws.Range("P97").Value = 18000
$ws.Range("Q97").Value = "$/caja 14 kilos empedrada"
$ws.Range("R97").Value = "Provincia de Limarí"
$ws.Range("S97").Value = 1286
$ws.Range("T97").Value = 14
